# Apply corrected matrix values to the "k_simp" worksheet.
# The diff only touches <v> (cached value) content of existing, purely
# numeric (non-formula) cells in rows 2-13 / columns A-L, so we just
# overwrite each changed cell's Value with its new number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A2"  = 771797499.4876
    "G2"  = -43041606.88665712
    "I2"  = -43041606.88665712

    "B3"  = 2408631579.893814
    "C3"  = 15105264.81417378
    "E3"  = -8631579.893813588
    "F3"  = 15105264.81417378

    "B4"  = 15105264.81417378
    "C4"  = 153660284.3409755
    "E4"  = -15105264.81417378
    "F4"  = 17291356.28194697
    "G4"  = 43041606.88665712
    "I4"  = 27041606.88665711

    "D5"  = 771797499.4876
    "J5"  = -43041606.88665712
    "L5"  = -43041606.88665712

    "B6"  = -8631579.893813588
    "C6"  = -15105264.81417378
    "E6"  = 2408631579.893814
    "F6"  = -15105264.81417378

    "B7"  = 15105264.81417378
    "C7"  = 17291356.28194697
    "E7"  = -15105264.81417378
    "F7"  = 153660284.3409755
    "J7"  = 43041606.88665712
    "L7"  = 27041606.88665711

    "A8"  = -43041606.88665712
    "C8"  = 43041606.88665712
    "G8"  = 728755892.6009429
    "I8"  = 43041606.88665712

    "H9"  = 1208631579.893814
    "I9"  = 15105264.81417378
    "K9"  = -8631579.893813588
    "L9"  = 15105264.81417378

    "A10" = -43041606.88665712
    "C10" = 27041606.88665711
    "G10" = 43041606.88665712
    "H10" = 15105264.81417378
    "I10" = 94618677.45431837
    "K10" = -15105264.81417378
    "L10" = 17291356.28194697

    "D11" = -43041606.88665712
    "F11" = 43041606.88665712
    "J11" = 728755892.6009429
    "L11" = 43041606.88665712

    "H12" = -8631579.893813588
    "I12" = -15105264.81417378
    "K12" = 1208631579.893814
    "L12" = -15105264.81417378

    "D13" = -43041606.88665712
    "F13" = 27041606.88665711
    "H13" = 15105264.81417378
    "I13" = 17291356.28194697
    "J13" = 43041606.88665712
    "K13" = -15105264.81417378
    "L13" = 94618677.45431837
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
